# Omaha_Cal_Info_GA05MOAS-GL493_00001.xlsx
# "added missing cal events, assigned OOI bar codes where necessary,
#  corrected integration events"

$wb = $excel.ActiveWorkbook

$moorings = $wb.Worksheets.Item("Moorings")
$calInfo  = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Corrected integration event for the GL493 controller (row 11) ---
# Previously E11 was blank and F11 held a bare "493" number (a stand-in /
# mistaken value). Replace with the real OOI barcode for the controller and
# a descriptive label, clearing the old numeric-cell formatting so the new
# text cells are written in the sheet's default (unstyled) format.
$calInfo.Range("E11").Value = "OL000108"
$calInfo.Range("F11").Value = "GL493 Controller"
$calInfo.Range("E11:F11").Style = "Normal"

# --- Update the active sheet / selection to reflect where work left off ---
$moorings.Range("C21").Select() | Out-Null

$calInfo.Select() | Out-Null
$calInfo.Range("E22").Select() | Out-Null
